# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which carry duplicate data tables.

$wb = $excel.ActiveWorkbook

# Row number (in sheet) -> new value for column F
$updates = @{
    2  = 197
    3  = 3231
    6  = 205
    7  = 1701
    19 = 27
    23 = 382
    24 = 227
    29 = 314
    30 = 2213
    34 = 440
    35 = 569
    38 = 347
    40 = 522
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
